# Actualización automática 2025-10-06 16:30:09
# Updates the "PRESUPUESTO" (budget) column G values on the "VENTA MENSUAL"
# sheet for a subset of advisor/client rows, and refreshes the total in G36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Row number (in sheet) -> new PRESUPUESTO value
$updates = @{
    2  = 1500
    4  = 500
    8  = 2000
    9  = 1500
    11 = 0
    12 = 1000
    13 = 1500
    15 = 3000
    19 = 2000
    20 = 4500
    23 = 750
    28 = 500
    29 = 10000
    32 = 1700
    33 = 0
    34 = 0
    35 = 1000
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}

# Recalculate the total row (G36) as the sum of G2:G35
$ws.Range("G36").Value = $ws.Application.WorksheetFunction.Sum($ws.Range("G2:G35"))
